$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 530.3418
$ws.Range("I15").Value = 530.3418
$ws.Range("K15").Value = 1591.0254
$ws.Range("M15").Value = -1422.0254
$ws.Range("H17").Value = 3006915
$ws.Range("J17").Value = 3006915
$ws.Range("L17").Value = 9020745
$ws.Range("N17").Value = -9021081
$ws.Range("H28").Value = 1555.5714
$ws.Range("J28").Value = 1983
$ws.Range("L28").Value = 1983
$ws.Range("N28").Value = -2953
$ws.Range("H33").Value = 11367346
$ws.Range("I33").Value = 16667558
$ws.Range("K33").Value = 16667558
$ws.Range("M33").Value = -16667329
$ws.Range("H41").Value = 163.78572
$ws.Range("I41").Value = 105.5
$ws.Range("J41").Value = 513.5
$ws.Range("K41").Value = 105.5
$ws.Range("L41").Value = 513.5
$ws.Range("M41").Value = 334.5
$ws.Range("N41").Value = -1393.5
$ws.Range("H55").Value = 253
$ws.Range("I55").Value = 204.4
$ws.Range("K55").Value = 204.4
$ws.Range("M55").Value = 9.599999999999994
$ws.Range("H64").Value = 7664.3335
$ws.Range("I64").Value = 7997
$ws.Range("K64").Value = 7997
$ws.Range("M64").Value = -7749
$ws.Range("H67").Value = 7664.3335
$ws.Range("I67").Value = 7997
$ws.Range("K67").Value = 7997
$ws.Range("M67").Value = -7139

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5011507
$ws.Range("I32").Value = 6582821
$ws.Range("K32").Value = 6582821
$ws.Range("M32").Value = -6582534
$ws.Range("H45").Value = 5375.8125
$ws.Range("I45").Value = 4999.923
$ws.Range("J45").Value = 7004.6665
$ws.Range("K45").Value = 4999.923
$ws.Range("L45").Value = 7004.6665
$ws.Range("M45").Value = -4622.923
$ws.Range("N45").Value = -7758.6665
$ws.Range("H61").Value = 23622.824
$ws.Range("I61").Value = 15320.125
$ws.Range("J61").Value = 31003
$ws.Range("K61").Value = 15320.125
$ws.Range("L61").Value = 31003
$ws.Range("M61").Value = -15108.125
$ws.Range("N61").Value = -31427
$ws.Range("H74").Value = 31887.428
$ws.Range("I74").Value = 5187.25
$ws.Range("J74").Value = 42567.5
$ws.Range("K74").Value = 5187.25
$ws.Range("L74").Value = 42567.5
$ws.Range("M74").Value = -4313.25
$ws.Range("N74").Value = -44315.5
$ws.Range("H77").Value = 31887.428
$ws.Range("I77").Value = 5187.25
$ws.Range("J77").Value = 42567.5
$ws.Range("K77").Value = 25936.25
$ws.Range("L77").Value = 212837.5
$ws.Range("M77").Value = -21568.25
$ws.Range("N77").Value = -221573.5
$ws.Range("H132").Value = 2574119
$ws.Range("I132").Value = 4264.92
$ws.Range("K132").Value = 12794.76
$ws.Range("M132").Value = -10264.76
$ws.Range("H136").Value = 23622.824
$ws.Range("I136").Value = 15320.125
$ws.Range("J136").Value = 31003
$ws.Range("K136").Value = 45960.375
$ws.Range("L136").Value = 93009
$ws.Range("M136").Value = -43410.375
$ws.Range("N136").Value = -98109

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5440.5
$ws.Range("I86").Value = 5017.9
$ws.Range("J86").Value = 6144.8335
$ws.Range("K86").Value = 5017.9
$ws.Range("L86").Value = 6144.8335
$ws.Range("M86").Value = -3894.9
$ws.Range("N86").Value = -8390.833500000001
$ws.Range("H89").Value = 5440.5
$ws.Range("I89").Value = 5017.9
$ws.Range("J89").Value = 6144.8335
$ws.Range("K89").Value = 25089.5
$ws.Range("L89").Value = 30724.1675
$ws.Range("M89").Value = -19473.5
$ws.Range("N89").Value = -41956.1675
$ws.Range("H105").Value = 3626.3
$ws.Range("I105").Value = 3626.3
$ws.Range("K105").Value = 3626.3
$ws.Range("M105").Value = -1879.3
$ws.Range("H107").Value = 9049.444
$ws.Range("I107").Value = 9555.625
$ws.Range("K107").Value = 9555.625
$ws.Range("M107").Value = -7635.625
$ws.Range("H134").Value = 11071.464
$ws.Range("I134").Value = 2831.9443
$ws.Range("J134").Value = 25902.6
$ws.Range("K134").Value = 8495.832900000001
$ws.Range("L134").Value = 77707.79999999999
$ws.Range("M134").Value = -5960.832900000001
$ws.Range("N134").Value = -82777.79999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4304
$ws.Range("I16").Value = 1092.875
$ws.Range("K16").Value = 1092.875
$ws.Range("M16").Value = -805.875
$ws.Range("H31").Value = 29857.947
$ws.Range("I31").Value = 18808.666
$ws.Range("K31").Value = 18808.666
$ws.Range("M31").Value = -18513.666
$ws.Range("H34").Value = 29857.947
$ws.Range("I34").Value = 18808.666
$ws.Range("K34").Value = 18808.666
$ws.Range("M34").Value = -18606.666
$ws.Range("H58").Value = 32556.23
$ws.Range("I58").Value = 17784.5
$ws.Range("K58").Value = 17784.5
$ws.Range("M58").Value = -17581.5
$ws.Range("H62").Value = 12327.2
$ws.Range("I62").Value = 13879.4
$ws.Range("J62").Value = 11551.1
$ws.Range("K62").Value = 13879.4
$ws.Range("L62").Value = 11551.1
$ws.Range("M62").Value = -13255.4
$ws.Range("N62").Value = -12799.1
$ws.Range("H65").Value = 12327.2
$ws.Range("I65").Value = 13879.4
$ws.Range("J65").Value = 11551.1
$ws.Range("K65").Value = 69397
$ws.Range("L65").Value = 57755.5
$ws.Range("M65").Value = -66277
$ws.Range("N65").Value = -63995.5
$ws.Range("H86").Value = 3557.9167
$ws.Range("I86").Value = 2666.818
$ws.Range("K86").Value = 2666.818
$ws.Range("M86").Value = -1543.818
$ws.Range("H89").Value = 3557.9167
$ws.Range("I89").Value = 2666.818
$ws.Range("K89").Value = 13334.09
$ws.Range("M89").Value = -7718.09
$ws.Range("H113").Value = 4304
$ws.Range("I113").Value = 1092.875
$ws.Range("K113").Value = 1092.875
$ws.Range("M113").Value = 1077.125
$ws.Range("H136").Value = 32556.23
$ws.Range("I136").Value = 17784.5
$ws.Range("K136").Value = 53353.5
$ws.Range("M136").Value = -50803.5
$ws.Range("H141").Value = 262890.62
$ws.Range("J141").Value = 262890.62
$ws.Range("L141").Value = 262890.62
$ws.Range("N141").Value = -273250.62

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 9764.799999999999
$ws.Range("J106").Value = 9764.799999999999
$ws.Range("L106").Value = 29294.4
$ws.Range("N106").Value = -31186.4
$ws.Range("H107").Value = 3907225.8
$ws.Range("J107").Value = 6251079
$ws.Range("L107").Value = 18753237
$ws.Range("N107").Value = -18757077
$ws.Range("H122").Value = 7423630.5
$ws.Range("I122").Value = 12458494
$ws.Range("K122").Value = 112126446
$ws.Range("M122").Value = -112123996
$ws.Range("H131").Value = 1488.79
$ws.Range("I131").Value = 1126.3334
$ws.Range("J131").Value = 1500
$ws.Range("K131").Value = 3379.0002
$ws.Range("L131").Value = 4500
$ws.Range("M131").Value = 1660.9998
$ws.Range("N131").Value = -14580

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 11188.567
$ws.Range("J80").Value = 12159.36
$ws.Range("L80").Value = 12159.36
$ws.Range("N80").Value = -14155.36
$ws.Range("H83").Value = 11188.567
$ws.Range("J83").Value = 12159.36
$ws.Range("L83").Value = 60796.8
$ws.Range("N83").Value = -70780.8
$ws.Range("H122").Value = 6189.6
$ws.Range("I122").Value = 1884.6
$ws.Range("J122").Value = 14799.6
$ws.Range("K122").Value = 5653.799999999999
$ws.Range("L122").Value = 44398.8
$ws.Range("M122").Value = -3203.799999999999
$ws.Range("N122").Value = -49298.8
$ws.Range("H132").Value = 7990.7744
$ws.Range("I132").Value = 3882.2273
$ws.Range("K132").Value = 11646.6819
$ws.Range("M132").Value = -9116.6819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 47823.72
$ws.Range("I122").Value = 67824.81
$ws.Range("K122").Value = 203474.43
$ws.Range("M122").Value = -201024.43
$ws.Range("H136").Value = 13828.167
$ws.Range("I136").Value = 14677.053
$ws.Range("K136").Value = 44031.159
$ws.Range("M136").Value = -41481.159

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2568.1765
$ws.Range("I113").Value = 1835.9
$ws.Range("K113").Value = 5507.700000000001
$ws.Range("M113").Value = -3337.700000000001
$ws.Range("H122").Value = 47625690
$ws.Range("I122").Value = 90912136
$ws.Range("J122").Value = 10599.5
$ws.Range("K122").Value = 272736408
$ws.Range("L122").Value = 31798.5
$ws.Range("M122").Value = -272733958
$ws.Range("N122").Value = -36698.5
